$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.867.66"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "3.190.52"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'537.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'144.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.07%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.25%  "
$ws.Range("D9").Value = "'7.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +3.97%  "
$ws.Range("D11").Value = "'0.431"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "3.747.33"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "'26.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "59.939.87"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").Value = "3.209.11"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'13.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "'8.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").Value = "'382.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").Value = "'70.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").Value = "'8.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.63%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "0.0₃0897"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").Value = "'1.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'22.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").Value = "'6.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("D33").Value = "'1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").Value = "'6.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.81%  "
$ws.Range("D35").Value = "'156.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("D36").Value = "'1.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "2.767.29"
$ws.Range("E37").Value = "  +5.22%  "
$ws.Range("D38").Value = "'25.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "'0.0714"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.86%  "
$ws.Range("D40").Value = "'1.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("E43").Value = "  +4.06%  "
$ws.Range("E44").Value = "  +5.41%  "
$ws.Range("D45").Value = "3.235.77"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").Value = "'6.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").Value = "  +6.39%  "
$ws.Range("D50").Value = "'20.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  -0.02%  "
